$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# ACCOUNT sheet (sheet1.xml) - add new row 16 describing ACC_SCHOOLNAME
# -----------------------------------------------------------------
$wsAccount = $wb.Worksheets.Item("ACCOUNT")
$wsAccount.Range("A16").Value = "ACC_SCHOOLNAME"
$wsAccount.Range("D16").Value = "CHAR"
$wsAccount.Range("E16").Value = 150
$wsAccount.Range("H16").Value = "REQUIRED"
$wsAccount.Activate() | Out-Null
$wsAccount.Range("A17").Select() | Out-Null

# -----------------------------------------------------------------
# STUDENT sheet (sheet2.xml) - bump size, insert ACC_FNAME/ACC_LNAME rows
# -----------------------------------------------------------------
$wsStudent = $wb.Worksheets.Item("STUDENT")
$wsStudent.Range("E5").Value = 100
$wsStudent.Range("A7").Value = "ACC_FNAME"
$wsStudent.Range("D7").Value = "CHAR"
$wsStudent.Range("E7").Value = 25
$wsStudent.Range("A8").Value = "ACC_LNAME"
$wsStudent.Range("D8").Value = "CHAR"
$wsStudent.Range("E8").Value = 25
$wsStudent.Activate() | Out-Null
$wsStudent.Range("E5").Select() | Out-Null

# -----------------------------------------------------------------
# INSTRUCTOR sheet (sheet3.xml) - bump size
# -----------------------------------------------------------------
$wsInstructor = $wb.Worksheets.Item("INSTRUCTOR")
$wsInstructor.Range("E5").Value = 100
$wsInstructor.Activate() | Out-Null
$wsInstructor.Range("E5").Select() | Out-Null

# -----------------------------------------------------------------
# CLASS sheet (sheet4.xml) - CLS_ID type CHAR -> INT
# -----------------------------------------------------------------
$wsClass = $wb.Worksheets.Item("CLASS")
$wsClass.Range("D5").Value = "INT"
$wsClass.Activate() | Out-Null
$wsClass.Range("D6").Select() | Out-Null

# -----------------------------------------------------------------
# COURSE sheet (sheet6.xml) - remove ACC_USERNAME row, shift Sdate/Edate rows up
# -----------------------------------------------------------------
$wsCourse = $wb.Worksheets.Item("COURSE")
$wsCourse.Range("A7").Value = "CLS_Sdate"
$wsCourse.Range("C7").Value = $null
$wsCourse.Range("D7").Value = "DATE"
$wsCourse.Range("H7").Value = "CLASS START DATE"
$wsCourse.Range("A8").Value = "CLS_Edate"
$wsCourse.Range("D8").Value = "DATE"
$wsCourse.Range("H8").Value = "CLASS END DATE"
$wsCourse.Range("A9:H9").Clear()
$wsCourse.Activate() | Out-Null
$wsCourse.Range("A7").Select() | Out-Null

# -----------------------------------------------------------------
# COMPLETION sheet (sheet7.xml) - add VARCHAR/CHAR field metadata
# -----------------------------------------------------------------
$wsCompletion = $wb.Worksheets.Item("COMPLETION")
$wsCompletion.Range("D6").Value = "VARCHAR"
$wsCompletion.Range("E6").Value = 100
$wsCompletion.Range("D7").Value = "CHAR"
$wsCompletion.Range("E7").Value = 5
$wsCompletion.Range("H7").Value = "FK to GRADE: GRD_CODE"
$wsCompletion.Activate() | Out-Null
$wsCompletion.Range("E8").Select() | Out-Null

$wsAccount.Activate() | Out-Null
